$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper cell (outside the table) used to stage plain-text values. Assigning a
# TEXT()-style formula ("="..."") forces the result to be a text value; copying
# that and pasting *values only* into the target cell stores it as shared-string
# text instead of Excel re-parsing it back into a number - matching how the
# other "X"/"Y" landmark coordinates in this sheet are stored as text.
$helper = $ws.Range("Z1")

function Paste-Text($targetAddress, $text) {
    $helper.Formula = "=""" + $text + """"
    $helper.Copy()
    $ws.Range($targetAddress).PasteSpecial(-4163)  # xlPasteValues
}

function Copy-Format($sourceAddress, $targetAddress) {
    $ws.Range($sourceAddress).Copy()
    $ws.Range($targetAddress).PasteSpecial(-4122)  # xlPasteFormats
}

# --- New landmark coordinates for LM 21 (row 22), LM 22 (row 23) and LM 23 (row 24) ---
# Row 22 already has a style applied (s="1"); only the values change.
Paste-Text "B22" "5.8"
Paste-Text "C22" "10.1"

# Rows 23 and 24 need both the new values and the matching style.
Paste-Text "B23" "5.9"
Paste-Text "C23" "4.0"
Paste-Text "D23" "0.2"

Paste-Text "B24" "1.4"
Paste-Text "C24" "7.2"
Paste-Text "D24" "0.2"

# Rows 25 and 26: column D gets the same "0.2" value; columns B and C stay
# empty but receive the same centred style as the rest of the table.
Paste-Text "D25" "0.2"
Paste-Text "D26" "0.2"

# Remove the helper cell so it leaves no trace in the sheet.
$helper.Clear()

# Apply the table's standard style (centered alignment) to every newly
# populated / still-blank cell in rows 23-26, columns B-D.
Copy-Format "B22" "B23"
Copy-Format "C22" "C23"
Copy-Format "D22" "D23"

Copy-Format "B22" "B24"
Copy-Format "C22" "C24"
Copy-Format "D22" "D24"

Copy-Format "B22" "B25"
Copy-Format "C22" "C25"
Copy-Format "D22" "D25"

Copy-Format "B22" "B26"
Copy-Format "C22" "C26"
Copy-Format "D22" "D26"

# Update the selected cell to where the user last clicked.
$ws.Range("B25").Select() | Out-Null
